# Updating Old File for core commit 9858844ccecc37046d166185cb936b938f965063
#
# The "Rules" sheet lists CodeQuality rules (Rule Key | Description | Type |
# Severity | Tags). This commit:
#   1. Adds a new Blocker-severity "CloudServiceIncompatibleWorkflowProcess"
#      Bug rule (inserted just above the existing "AEM Rules:AEM-3" row).
#   2. Removes the deprecated "CQRules:CQBP-84--dependencies" rule row.
#   3. Removes the old Major-severity "CloudServiceIncompatibleWorkflowProcess"
#      Bug rule row (superseded by the new Blocker-severity row added above).
#   4. Adds a new Minor-severity "IndexDamAssetLucene" Bug rule (inserted
#      just above the existing "ClientlibProxyResource" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new Blocker "CloudServiceIncompatibleWorkflowProcess" rule
#    above row 36 ("AEM Rules:AEM-3").
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "CloudServiceIncompatibleWorkflowProcess"
$ws.Range("B36").Value = "Usage of Cloud Service Incompatible Workflow Processes"
$ws.Range("C36").Value = "Bug"
$ws.Range("D36").Value = "Blocker"
$ws.Range("E36").Value = "aem,cloud-service-compatibility"

# 2) Delete the "CQRules:CQBP-84--dependencies" rule row (now row 39, after
#    the insert above shifted everything below row 36 down by one).
$ws.Rows.Item(39).Delete()

# 3) Delete the old Major "CloudServiceIncompatibleWorkflowProcess" rule row
#    (row 45).
$ws.Rows.Item(45).Delete()

# 4) Insert the new Minor "IndexDamAssetLucene" rule above row 70
#    ("ClientlibProxyResource").
$ws.Rows.Item(70).Insert()
$ws.Range("A70").Value = "IndexDamAssetLucene"
$ws.Range("B70").Value = "Index customizations of the damAssetLucene Oak index should be properly structured."
$ws.Range("C70").Value = "Bug"
$ws.Range("D70").Value = "Minor"
$ws.Range("E70").Value = "aem,cloud-service-compatibility"

# Restore the cell selection recorded in the saved view state.
$ws.Range("E70").Select()
